$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.06958889789029854
    "D2" = 0.01865746595503026
    "E2" = 0.1253897863139528
    "F2" = 0.8338999087808929
    "G2" = 0.6853892590213491
    "H2" = 0.7456160989791414
    "I2" = 0.4797590421002109
    "K2" = 0.6878161537638903
    "M2" = 0.3028016575899173
    "N2" = 1.398722012147433
    "B3" = 0.06157872697281164
    "D3" = 0.01826407897105753
    "E3" = 0.1182100319822865
    "F3" = 0.819160508377351
    "G3" = 0.6708967082266781
    "H3" = 0.7438408254652558
    "I3" = 0.485014859274127
    "K3" = 0.6016167525981757
    "M3" = 0.271200669196098
    "N3" = 1.415878419800123
    "B4" = 0.05666112932409817
    "D4" = 0.01802051241072178
    "E4" = 0.1138949797749618
    "F4" = 0.8107056705196953
    "G4" = 0.6625382463223701
    "H4" = 0.7431806631114597
    "I4" = 0.4884748238995886
    "K4" = 0.5486881255251035
    "M4" = 0.2519036850950584
    "N4" = 1.426943602025025
    "B5" = 0.05465750990107665
    "D5" = 0.01792075682837435
    "E5" = 0.112159804951169
    "F5" = 0.8074095121924074
    "G5" = 0.6592673441217016
    "H5" = 0.7430195921829608
    "I5" = 0.48994331373782
    "K5" = 0.5271188101892506
    "M5" = 0.2440664085499407
    "N5" = 1.431586183321083
    "B6" = 0.05432483607913241
    "D6" = 0.01790416254644711
    "E6" = 0.1118730777118557
    "F6" = 0.8068711902939611
    "G6" = 0.6587323652439068
    "H6" = 0.7429993623225783
    "I6" = 0.4901906892804622
    "K6" = 0.5235372134967804
    "M6" = 0.2427666208832377
    "N6" = 1.432365135936329
    "B7" = 0.05663410620226728
    "D7" = 0.01801916908512169
    "E7" = 0.1138714847872677
    "F7" = 0.8106606135923329
    "G7" = 0.6624935869499353
    "H7" = 0.7431780539564556
    "I7" = 0.4884943915394722
    "K7" = 0.5483972360989355
    "M7" = 0.2517978823269686
    "N7" = 1.427005673404151
    "B8" = 0.06682694961325808
    "D8" = 0.01852225175682065
    "E8" = 0.1228946928539187
    "F8" = 0.8286939699133598
    "G8" = 0.68027980525099
    "H8" = 0.7449146781294473
    "I8" = 0.4815228968491923
    "K8" = 0.6580949522713695
    "M8" = 0.2918834196355604
    "N8" = 1.40452724888517
    "B9" = 0.08681396115937901
    "D9" = 0.01949237789740366
    "E9" = 0.1413411631070502
    "F9" = 0.8688018985287016
    "G9" = 0.7194712097754064
    "H9" = 0.7517386422644563
    "I9" = 0.4697008710587802
    "K9" = 0.8732085228993753
    "M9" = 0.3713505521204041
    "N9" = 1.364665417619824
    "B10" = 0.1014904173427311
    "D10" = 0.02019472959581847
    "E10" = 0.1553693624735359
    "F10" = 0.9011963154165841
    "G10" = 0.7509369080692352
    "H10" = 0.7588486718985621
    "I10" = 0.4621441648276878
    "K10" = 1.03128208921629
    "M10" = 0.4302888706820482
    "N10" = 1.337955575901201
    "B11" = 0.1081639210876517
    "D11" = 0.0205119137402221
    "E11" = 0.1618582172864222
    "F11" = 0.9165767556580988
    "G11" = 0.7658410851003907
    "H11" = 0.7625412783351351
    "I11" = 0.4589519337000176
    "K11" = 1.103208110528556
    "M11" = 0.4572285004743719
    "N11" = 1.326365178151516
    "B12" = 0.1106904339581689
    "D12" = 0.02063168211895317
    "E12" = 0.1643310936420903
    "F12" = 0.9224940853589487
    "G12" = 0.7715704757697495
    "H12" = 0.7640056645557252
    "I12" = 0.4577784382946142
    "K12" = 1.130447432383505
    "M12" = 0.4674487081792478
    "N12" = 1.322056863995835
    "B13" = 0.1101463337117394
    "D13" = 0.02060590322140143
    "E13" = 0.1637978137511453
    "F13" = 0.9212155346832134
    "G13" = 0.7703327363062158
    "H13" = 0.763687341048211
    "I13" = 0.4580295992418684
    "K13" = 1.124580846292929
    "M13" = 0.4652467645655349
    "N13" = 1.322981143606135
    "B14" = 0.1083717918541254
    "D14" = 0.02052177406901379
    "E14" = 0.162061346465542
    "F14" = 0.9170617095190892
    "G14" = 0.7663107283538011
    "H14" = 0.7626604289351633
    "I14" = 0.4588546810978578
    "K14" = 1.105449055990505
    "M14" = 0.4580689454790274
    "N14" = 1.326009111776301
    "B15" = 0.1072847504195522
    "D15" = 0.02047019768600578
    "E15" = 0.1609997605889362
    "F15" = 0.9145295090121124
    "G15" = 0.7638582857756262
    "H15" = 0.7620400262304372
    "I15" = 0.4593646705338408
    "K15" = 1.093730606694749
    "M15" = 0.453674775570363
    "N15" = 1.327874347103792
    "B16" = 0.101054214404499
    "D16" = 0.02017395343421668
    "E16" = 0.1549474821831396
    "F16" = 0.9002041681927437
    "G16" = 0.74997481593212
    "H16" = 0.7586165876756468
    "I16" = 0.4623577261835763
    "K16" = 1.026581905553996
    "M16" = 0.4285309102096164
    "N16" = 1.338724320830661
    "B17" = 0.09723111209878255
    "D17" = 0.01999161680616268
    "E17" = 0.1512622794485594
    "F17" = 0.8915813320514019
    "G17" = 0.7416093814059508
    "H17" = 0.7566339106832629
    "I17" = 0.4642567463271874
    "K17" = 0.9853928297657717
    "M17" = 0.4131390390733287
    "N17" = 1.345523998809208
    "B18" = 0.09503190364840464
    "D18" = 0.0198865238279744
    "E18" = 0.1491527381650002
    "F18" = 0.886682297863274
    "G18" = 0.7368533266566146
    "H18" = 0.7555366457704622
    "I18" = 0.465372103094456
    "K18" = 0.9617036171944733
    "M18" = 0.404298098972717
    "N18" = 1.34948769637311
    "B19" = 0.09428725018126727
    "D19" = 0.01985090402169121
    "E19" = 0.1484402090805261
    "F19" = 0.885033961530965
    "G19" = 0.7352525252919691
    "H19" = 0.7551725307642982
    "I19" = 0.4657537079924232
    "K19" = 0.9536831419092096
    "M19" = 0.4013067700928445
    "N19" = 1.350838782441432
    "B20" = 0.09763811621589014
    "D20" = 0.02001104945521348
    "E20" = 0.1516535295033421
    "F20" = 0.8924929730778217
    "G20" = 0.7424941453838301
    "H20" = 0.756840506097177
    "I20" = 0.4640522022236127
    "K20" = 0.989777306190831
    "M20" = 0.4147762799648547
    "N20" = 1.344794705624858
    "B21" = 0.1088930355000315
    "D21" = 0.02054649420216847
    "E21" = 0.1625709612342092
    "F21" = 0.9182792583943069
    "G21" = 0.7674897641068412
    "H21" = 0.762960263137245
    "I21" = 0.4586113752074894
    "K21" = 1.111068458940736
    "M21" = 0.4601767334892486
    "N21" = 1.325117531975753
    "B22" = 0.1162452049927793
    "D22" = 0.02089443869765262
    "E22" = 0.1697976921437387
    "F22" = 0.935674913795765
    "G22" = 0.7843245093790188
    "H22" = 0.7673451089029584
    "I22" = 0.4552614536973749
    "K22" = 1.190353793963368
    "M22" = 0.4899580145983577
    "N22" = 1.312727935574669
    "B23" = 0.1123216006879488
    "D23" = 0.02070892001968971
    "E23" = 0.1659321891893484
    "F23" = 0.926340698345868
    "G23" = 0.7752936550378422
    "H23" = 0.7649695238224297
    "I23" = 0.4570305060515203
    "K23" = 1.148036394183009
    "M23" = 0.4740530668160545
    "N23" = 1.319297376902389
    "B24" = 0.09745411355333999
    "D24" = 0.0200022647797482
    "E24" = 0.1514766169339197
    "F24" = 0.8920806384466289
    "G24" = 0.7420939774667943
    "H24" = 0.7567469716120456
    "I24" = 0.4641446031493537
    "K24" = 0.987795112641578
    "M24" = 0.4140360583081559
    "N24" = 1.345124249398722
    "B25" = 0.08140782391967605
    "D25" = 0.01923173187576666
    "E25" = 0.1362686321030893
    "F25" = 0.8574401253461872
    "G25" = 0.7084027967966051
    "H25" = 0.7495253284996437
    "I25" = 0.4727009136054434
    "K25" = 0.8150117062861
    "M25" = 0.3497574968171264
    "N25" = 1.374996993277753
}

foreach ($cellRef in $values.Keys) {
    $ws.Range($cellRef).Value = $values[$cellRef]
}
